$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text that previously said "aaaa" to the new text
$ws.Range("B2").Value = "to be continue…."

# Clear the cells that were removed from row 2 (C2:E2)
$ws.Range("C2:E2").ClearContents()

# Update the active selection to C7
$ws.Range("C7").Select()
